$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.520080089569092
$ws.Range("B1").Value = 2.776425361633301
$ws.Range("C1").Value = 1.781080365180969
$ws.Range("D1").Value = 1.089894413948059
$ws.Range("E1").Value = 0.5591490864753723
